$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.965.46"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "1.643.67"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5210"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06363"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.68"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07686"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.423"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.630.36"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").Value = "1.867.52"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5510"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").Value = "0.0₅8256"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "25.983.58"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.696"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.16"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.23"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1238"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.392"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.91"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05905"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.77%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.389"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.392"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.645"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9920"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.392"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.749"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5626"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.79%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8532"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "1.033.13"
$ws.Range("E43").Value = "  -6.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.70"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.87%  "
$ws.Range("D45").Value = "1.791.52"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.59"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.044"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4218"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.884"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.24%  "
